# Replace the RNA-style "U" in the Codon column (column A) with the
# DNA-style "T", e.g. "UUU" -> "TTT", leaving the Amino (B) and
# Frequency (C) columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val.ToString().Replace("U", "T")
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Mirror the resulting full-column selection seen in the saved file.
$ws.Columns("A").Select()
